$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking price strings (e.g. "326.70", "1.001") that
# must stay as literal text (matching the source "inlineStr" cells, which used
# thousand-dot / no-thousand-separator formatting rather than real numbers).
# Excel auto-converts a bare numeric-looking assignment to a Number, which loses
# trailing zeros / precision, so we force Text format, assign, then clear the
# explicit formatting again (ClearFormats) so no stray style is left behind on
# the cell - only the literal text survives.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '27.711.84'
$ws.Range('E2').Value = '  +0.78%  '
Set-TextValue $ws.Range('D3') '1.775.58'
$ws.Range('E3').Value = '  +1.47%  '
$ws.Range('E4').Value = '  +0.02%  '
Set-TextValue $ws.Range('D5') '326.70'
$ws.Range('E5').Value = '  +0.79%  '
Set-TextValue $ws.Range('D6') '1.001'
$ws.Range('E6').Value = '  +0.03%  '
Set-TextValue $ws.Range('D7') '0.4619'
$ws.Range('E7').Value = '  +3.63%  '
Set-TextValue $ws.Range('D8') '0.3582'
$ws.Range('E8').Value = '  -0.62%  '
$ws.Range('E9').Value = '  -0.31%  '
Set-TextValue $ws.Range('D10') '41.76'
$ws.Range('E10').Value = '  -0.61%  '
$ws.Range('E11').Value = '  +0.81%  '
Set-TextValue $ws.Range('D12') '1.000'
$ws.Range('E12').Value = '  +0.06%  '
$ws.Range('E13').Value = '  +0.88%  '
Set-TextValue $ws.Range('D14') '6.034'
$ws.Range('E14').Value = '  +0.23%  '
Set-TextValue $ws.Range('D15') '7.237'
$ws.Range('E15').Value = '  +1.51%  '
Set-TextValue $ws.Range('D16') '1.781.07'
$ws.Range('E16').Value = '  +1.77%  '
Set-TextValue $ws.Range('D17') '93.53'
$ws.Range('E17').Value = '  +1.03%  '
$ws.Range('E18').Value = '  -0.28%  '
Set-TextValue $ws.Range('D19') '0.06406'
$ws.Range('E19').Value = '  +0.04%  '
$ws.Range('E20').Value = '  +0.01%  '
Set-TextValue $ws.Range('D21') '17.06'
$ws.Range('E21').Value = '  +1.42%  '
Set-TextValue $ws.Range('D22') '5.781'
$ws.Range('E22').Value = '  -1.14%  '
Set-TextValue $ws.Range('D23') '27.795.81'
$ws.Range('E23').Value = '  +0.92%  '
Set-TextValue $ws.Range('D24') '11.27'
$ws.Range('E24').Value = '  +1.02%  '
Set-TextValue $ws.Range('D25') '2.082'
$ws.Range('E25').Value = '  -0.62%  '
Set-TextValue $ws.Range('D26') '164.55'
$ws.Range('E26').Value = '  +1.59%  '
Set-TextValue $ws.Range('D27') '20.26'
$ws.Range('E27').Value = '  -1.06%  '
Set-TextValue $ws.Range('D28') '1.979.02'
$ws.Range('E28').Value = '  +1.45%  '
Set-TextValue $ws.Range('D29') '2.162'
$ws.Range('E29').Value = '  +3.87%  '
Set-TextValue $ws.Range('D30') '125.81'
$ws.Range('E30').Value = '  +0.81%  '
Set-TextValue $ws.Range('D31') '1.088'
$ws.Range('E31').Value = '  +0.45%  '
Set-TextValue $ws.Range('D32') '0.09231'
$ws.Range('E32').Value = '  +2.46%  '
Set-TextValue $ws.Range('D33') '3.671'
$ws.Range('E33').Value = '  +0.40%  '
Set-TextValue $ws.Range('D34') '5.528'
$ws.Range('E34').Value = '  +0.14%  '
Set-TextValue $ws.Range('D35') '11.79'
$ws.Range('E35').Value = '  -1.82%  '
Set-TextValue $ws.Range('D36') '0.02292'
$ws.Range('E36').Value = '  -0.33%  '
Set-TextValue $ws.Range('D37') '0.06169'
$ws.Range('E37').Value = '  +2.77%  '
Set-TextValue $ws.Range('D38') '0.2085'
$ws.Range('E38').Value = '  +0.05%  '
Set-TextValue $ws.Range('D39') '0.6310'
$ws.Range('E39').Value = '  -0.58%  '
Set-TextValue $ws.Range('D40') '4.949'
$ws.Range('E40').Value = '  +0.19%  '
Set-TextValue $ws.Range('D41') '1.181'
$ws.Range('E41').Value = '  -1.76%  '
Set-TextValue $ws.Range('D42') '1.391'
$ws.Range('E42').Value = '  +0.44%  '
$ws.Range('E43').Value = '  -0.07%  '
Set-TextValue $ws.Range('D44') '13.17'
$ws.Range('E44').Value = '  +0.42%  '
Set-TextValue $ws.Range('D45') '3.735'
$ws.Range('E45').Value = '  +0.73%  '
Set-TextValue $ws.Range('D46') '0.5880'
$ws.Range('E46').Value = '  -0.04%  '
Set-TextValue $ws.Range('D47') '122.23'
$ws.Range('E47').Value = '  +0.75%  '
Set-TextValue $ws.Range('D48') '1.944'
$ws.Range('E48').Value = '  -0.32%  '
Set-TextValue $ws.Range('D49') '0.06934'
$ws.Range('E49').Value = '  +1.09%  '
Set-TextValue $ws.Range('D50') '1.134'
$ws.Range('E50').Value = '  -1.38%  '
Set-TextValue $ws.Range('D51') '72.18'
$ws.Range('E51').Value = '  +0.33%  '
